$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price (D) cells that will receive numeric-looking text
# as Text ("@") so Excel stores them as strings, matching the source data
# (prices like "209.82" must stay literal text, not be parsed as numbers).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values
$ws.Range("D2").Value = '80.536.56'
$ws.Range("E2").Value = '  +5.27%  '
$ws.Range("D3").Value = '3.177.92'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '209.82'
$ws.Range("E5").Value = '  +5.04%  '
$ws.Range("D6").Value = '626.10'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = '0.275'
$ws.Range("E7").Value = '  +28.18%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.589'
$ws.Range("E9").Value = '  +6.49%  '
$ws.Range("D10").Value = '3.174.55'
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("D11").Value = '0.591'
$ws.Range("E11").Value = '  +26.55%  '
$ws.Range("D12").Value = '0.0000258'
$ws.Range("E12").Value = '  +28.81%  '
$ws.Range("E13").Value = '  +1.53%  '
$ws.Range("D14").Value = '3.757.61'
$ws.Range("E14").Value = '  +1.99%  '
$ws.Range("D15").Value = '5.27'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '32.03'
$ws.Range("E16").Value = '  +7.79%  '
$ws.Range("D17").Value = '80.382.02'
$ws.Range("E17").Value = '  +5.16%  '
$ws.Range("D18").Value = '3.170.63'
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("D19").Value = '14.20'
$ws.Range("E19").Value = '  +3.41%  '
$ws.Range("D20").Value = '3.02'
$ws.Range("E20").Value = '  +8.55%  '
$ws.Range("D21").Value = '9.18'
$ws.Range("E21").Value = '  -0.51%  '
$ws.Range("D22").Value = '438.09'
$ws.Range("E22").Value = '  +13.12%  '
$ws.Range("D23").Value = '5.20'
$ws.Range("E23").Value = '  +14.30%  '
$ws.Range("D24").Value = '6.94'
$ws.Range("E24").Value = '  +7.35%  '
$ws.Range("E25").Value = '  +2.21%  '
$ws.Range("D26").Value = '75.89'
$ws.Range("E26").Value = '  +4.18%  '
$ws.Range("D27").Value = '4.68'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '10.89'
$ws.Range("E28").Value = '  +4.31%  '
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '0.0000121'
$ws.Range("E30").Value = '  +8.28%  '
$ws.Range("D31").Value = '0.996'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("D32").Value = '8.91'
$ws.Range("E32").Value = '  +5.40%  '
$ws.Range("D33").Value = '560.77'
$ws.Range("E33").Value = '  +8.95%  '
$ws.Range("D34").Value = '1.46'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '0.151'
$ws.Range("E35").Value = '  +13.20%  '
$ws.Range("D36").Value = '1.99'
$ws.Range("E36").Value = '  +1.98%  '
$ws.Range("D37").Value = '22.96'
$ws.Range("E37").Value = '  +8.72%  '
$ws.Range("E38").Value = '  +18.80%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '0.406'
$ws.Range("E40").Value = '  +6.04%  '
$ws.Range("D41").Value = '20.78'
$ws.Range("E41").Value = '  +3.57%  '
$ws.Range("D42").Value = '162.71'
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '5.62'
$ws.Range("E44").Value = '  +4.54%  '
$ws.Range("D45").Value = '189.70'
$ws.Range("E45").Value = '  -4.41%  '
$ws.Range("D46").Value = '1.81'
$ws.Range("E46").Value = '  +6.73%  '
$ws.Range("D47").Value = '2.70'
$ws.Range("E47").Value = '  +7.79%  '
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("D49").Value = '1.30'
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").Value = '42.76'
$ws.Range("E50").Value = '  +3.73%  '
$ws.Range("D51").Value = '4.24'
$ws.Range("E51").Value = '  +6.99%  '
